$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The municipio-nombre metadata column (D) is re-classified from a measure
# to a curated dimension, per the re-processed data.
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"
